$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends with:
#   row 162: last data row (2023-05-31, 1991.197929)
#   row 164: "数据来源：Wind" footer
#
# We need to insert a new data row (2023-06-30, 1195.261165) directly after
# row 162, pushing the footer row down from 164 to 165.

# Shift the footer row (and everything from it down) one row lower, opening
# up row 163 for the new data row.
$ws.Rows.Item(164).Insert()

# Copy the formatting (date format / number format) from the previous data
# row so the new row matches the rest of the table exactly.
$ws.Range("A162:B162").Copy()
$ws.Range("A163:B163").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new data values.
$ws.Cells.Item(163, 1).Value = 45107.0
$ws.Cells.Item(163, 2).Value = 1195.261165
